$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "DO"
$ws.Range("C2").Value = "M1"
$ws.Range("E2").Value = "M3"
$ws.Range("H2").Value = "M1"
$ws.Range("I2").Value = "DO"
$ws.Range("J2").Value = "M3"
$ws.Range("K2").Value = "M1"
$ws.Range("O2").Value = "M1"
$ws.Range("P2").Value = "M1"
$ws.Range("Q2").Value = "DO"
$ws.Range("R2").Value = "M3"
$ws.Range("U2").Value = "M1"
$ws.Range("W2").Value = "A1"
$ws.Range("X2").Value = "A1"
$ws.Range("Z2").Value = "DO"
$ws.Range("B3").Value = "M1"
$ws.Range("C3").Value = "M1"
$ws.Range("D3").Value = "A1"
$ws.Range("E3").Value = "DO"
$ws.Range("F3").Value = "M2"
$ws.Range("G3").Value = "A1"
$ws.Range("J3").Value = "A2"
$ws.Range("K3").Value = "DO"
$ws.Range("L3").Value = "M2"
$ws.Range("M3").Value = "A1"
$ws.Range("N3").Value = "M2"
$ws.Range("P3").Value = "A2"
$ws.Range("Q3").Value = "M2"
$ws.Range("R3").Value = "M1"
$ws.Range("T3").Value = "A2"
$ws.Range("U3").Value = "M2"
$ws.Range("V3").Value = "DO"
$ws.Range("X3").Value = "DO"
$ws.Range("AA3").Value = "A1"
$ws.Range("AB3").Value = "M2"
$ws.Range("AC3").Value = "A2"
$ws.Range("B4").Value = "M3"
$ws.Range("C4").Value = "DO"
$ws.Range("E4").Value = "A1"
$ws.Range("F4").Value = "M1"
$ws.Range("H4").Value = "A1"
$ws.Range("K4").Value = "M3"
$ws.Range("M4").Value = "M1"
$ws.Range("O4").Value = "M3"
$ws.Range("Q4").Value = "A1"
$ws.Range("S4").Value = "M1"
$ws.Range("T4").Value = "DO"
$ws.Range("U4").Value = "M3"
$ws.Range("W4").Value = "M3"
$ws.Range("X4").Value = "M1"
$ws.Range("AA4").Value = "M1"
$ws.Range("AB4").Value = "A1"
$ws.Range("B5").Value = "A2"
$ws.Range("C5").Value = "A1"
$ws.Range("D5").Value = "M2"
$ws.Range("F5").Value = "A1"
$ws.Range("H5").Value = "DO"
$ws.Range("I5").Value = "A2"
$ws.Range("J5").Value = "M1"
$ws.Range("K5").Value = "A1"
$ws.Range("L5").Value = "A2"
$ws.Range("M5").Value = "DO"
$ws.Range("P5").Value = "DO"
$ws.Range("Q5").Value = "M1"
$ws.Range("S5").Value = "M1"
$ws.Range("T5").Value = "M2"
$ws.Range("U5").Value = "A1"
$ws.Range("W5").Value = "M2"
$ws.Range("X5").Value = "M2"
$ws.Range("Z5").Value = "M2"
$ws.Range("AA5").Value = "DO"
$ws.Range("AB5").Value = "M1"
$ws.Range("AC5").Value = "M1"
$ws.Range("D6").Value = "A2"
$ws.Range("E6").Value = "M2"
$ws.Range("F6").Value = "M2"
$ws.Range("H6").Value = "A1"
$ws.Range("I6").Value = "DO"
$ws.Range("K6").Value = "A1"
$ws.Range("L6").Value = "A1"
$ws.Range("M6").Value = "M2"
$ws.Range("P6").Value = "DO"
$ws.Range("Q6").Value = "M1"
$ws.Range("S6").Value = "M2"
$ws.Range("T6").Value = "M1"
$ws.Range("W6").Value = "M1"
$ws.Range("Y6").Value = "A1"
$ws.Range("AC6").Value = "DO"
$ws.Range("D7").Value = "DO"
$ws.Range("E7").Value = "A1"
$ws.Range("G7").Value = "M1"
$ws.Range("H7").Value = "M3"
$ws.Range("K7").Value = "M3"
$ws.Range("O7").Value = "M1"
$ws.Range("S7").Value = "A1"
$ws.Range("U7").Value = "M1"
$ws.Range("Y7").Value = "M3"
$ws.Range("Z7").Value = "M1"
$ws.Range("AA7").Value = "DO"
$ws.Range("AC7").Value = "A1"
$ws.Range("B8").Value = "M1"
$ws.Range("C8").Value = "M2"
$ws.Range("D8").Value = "A2"
$ws.Range("F8").Value = "DO"
$ws.Range("G8").Value = "M1"
$ws.Range("H8").Value = "M2"
$ws.Range("I8").Value = "M1"
$ws.Range("K8").Value = "M3"
$ws.Range("L8").Value = "A1"
$ws.Range("M8").Value = "A1"
$ws.Range("N8").Value = "M1"
$ws.Range("P8").Value = "M2"
$ws.Range("S8").Value = "A1"
$ws.Range("U8").Value = "M1"
$ws.Range("V8").Value = "A1"
$ws.Range("W8").Value = "M1"
$ws.Range("X8").Value = "A2"
$ws.Range("Y8").Value = "M1"
$ws.Range("Z8").Value = "A1"
$ws.Range("AB8").Value = "M2"
$ws.Range("AC8").Value = "DO"
$ws.Range("B9").Value = "A2"
$ws.Range("C9").Value = "DO"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = "M1"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "A1"
$ws.Range("J9").Value = "M3"
$ws.Range("K9").Value = "A2"
$ws.Range("L9").Value = "M1"
$ws.Range("M9").Value = "M3"
$ws.Range("N9").Value = "A2"
$ws.Range("P9").Value = "M1"
$ws.Range("Q9").Value = "A2"
$ws.Range("R9").Value = "M2"
$ws.Range("S9").Value = "DO"
$ws.Range("U9").Value = "A1"
$ws.Range("V9").Value = "M2"
$ws.Range("W9").Value = "A2"
$ws.Range("X9").Value = "DO"
$ws.Range("AA9").Value = "A2"
$ws.Range("AB9").Value = "M1"
$ws.Range("AC9").Value = "A1"
$ws.Range("B10").Value = "M2"
$ws.Range("C10").Value = "A1"
$ws.Range("E10").Value = "M2"
$ws.Range("F10").Value = "A2"
$ws.Range("G10").Value = "M1"
$ws.Range("H10").Value = "DO"
$ws.Range("I10").Value = "A1"
$ws.Range("M10").Value = "DO"
$ws.Range("N10").Value = "M2"
$ws.Range("O10").Value = "M1"
$ws.Range("P10").Value = "A2"
$ws.Range("Q10").Value = "M2"
$ws.Range("T10").Value = "A1"
$ws.Range("U10").Value = "A2"
$ws.Range("V10").Value = "DO"
$ws.Range("X10").Value = "M1"
$ws.Range("Z10").Value = "A2"
$ws.Range("AA10").Value = "M2"
